$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update time-of-day fraction on a handful of early datetime cells ---
# (09:59:59.999... -> 11:00:00.000...; value represents day-fraction shift)
$ws.Range("A6").Value = 37347.45833333334
$ws.Range("A18").Value = 37712.45833333334
$ws.Range("A30").Value = 38078.45833333334
$ws.Range("A42").Value = 38443.45833333334
$ws.Range("A60").Value = 38991.45833333334

# --- Refresh OHLC figures for the last few existing observations ---
$ws.Range("C252:F252").Value = 1703600000000
$ws.Range("C253:F253").Value = 1719900000000
$ws.Range("C254:F254").Value = 1749500000000

# --- Append the new latest observation (row 257) ---
$ws.Range("A257").Value = 45047.41666666666
$ws.Range("B257").Value = "ECONOMICS:AEM2"
$ws.Range("C257:F257").Value = 1855323000000
$ws.Range("G257").Value = 0

# Match the date-formatted style used by the rest of column A (border,
# centered/top alignment, custom date/time number format) on the new cell.
$ws.Range("A256").Copy()
$ws.Range("A257").PasteSpecial(-4122)
$excel.CutCopyMode = 0
